# New translations: duplicate the existing "french" (column E) translation
# text into a brand-new column F (rows 2-6), using the same styling that
# the bulk of column E already uses (the style found on E3:E6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the cell format from E3 (the shared style used by E3:E6) onto the
#    new F2:F6 range so the new cells pick up the identical style, without
#    introducing any new style/font/fill definitions.
$ws.Range("E3").Copy()
$ws.Range("F2:F6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 2) Populate F2:F6 with the same text currently shown in E2:E6.
for ($row = 2; $row -le 6; $row++) {
    $srcCell = $ws.Cells.Item($row, 5)   # column E
    $dstCell = $ws.Cells.Item($row, 6)   # column F
    $dstCell.Value2 = $srcCell.Value2
}
